# Insert a new weekly price record as row 190 in the "Uva" (grape) sheet.
# This pushes the previous rows 190-211 down to 191-212 (dimension grows
# from A1:T211 to A1:T212), matching the upstream commit
# "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 190..211 down by one, creating a fresh blank row 190.
$ws.Rows.Item(190).EntireRow.Insert()

# Populate the new row 190 with the new weekly record.
$ws.Range("A190").Value = 2
$ws.Range("B190").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C190").Value = "Coquimbo"
$ws.Range("D190").Value = 45008
$ws.Range("E190").Value = 4
$ws.Range("F190").Value = "Fruta"
$ws.Range("G190").Value = 100109
$ws.Range("H190").Value = "Uva"
$ws.Range("I190").Value = 100109001
$ws.Range("J190").Value = "Uva"
$ws.Range("K190").Value = "Moscatel rosada"
$ws.Range("L190").Value = "Primera"
$ws.Range("M190").Value = 500
$ws.Range("N190").Value = 13000
$ws.Range("O190").Value = 14000
$ws.Range("P190").Value = 13500
$ws.Range("Q190").Value = "`$/bandeja 18 kilos"
$ws.Range("R190").Value = "Provincia de Limarí"
$ws.Range("S190").Value = 750
$ws.Range("T190").Value = 18
